$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 265
$ws.Range("I33").Value = 289.33334
$ws.Range("J33").Value = 216.33333
$ws.Range("K33").Value = 289.33334
$ws.Range("L33").Value = 216.33333
$ws.Range("M33").Value = -60.33334000000002
$ws.Range("N33").Value = -674.3333299999999
$ws.Range("H100").Value = 1242.8695
$ws.Range("I100").Value = 981.625
$ws.Range("K100").Value = 981.625
$ws.Range("M100").Value = -440.625
$ws.Range("H132").Value = 1731.1818
$ws.Range("I132").Value = 1472.8572
$ws.Range("J132").Value = 2183.25
$ws.Range("K132").Value = 4418.571599999999
$ws.Range("L132").Value = 6549.75
$ws.Range("M132").Value = -1888.571599999999
$ws.Range("N132").Value = -11609.75
$ws.Range("H134").Value = 84069.75
$ws.Range("J134").Value = 84069.75
$ws.Range("L134").Value = 84069.75
$ws.Range("N134").Value = -94209.75
$ws.Range("H138").Value = 1888.8292
$ws.Range("I138").Value = 1324.2858
$ws.Range("J138").Value = 2309.2341
$ws.Range("K138").Value = 3972.8574
$ws.Range("L138").Value = 6927.702300000001
$ws.Range("M138").Value = 1167.1426
$ws.Range("N138").Value = -17207.7023

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6267.2856
$ws.Range("I61").Value = 5313.121
$ws.Range("J61").Value = 9765.888999999999
$ws.Range("K61").Value = 5313.121
$ws.Range("L61").Value = 9765.888999999999
$ws.Range("M61").Value = -5101.121
$ws.Range("N61").Value = -10189.889
$ws.Range("H74").Value = 3981.8108
$ws.Range("I74").Value = 2743.8572
$ws.Range("J74").Value = 7833.222
$ws.Range("K74").Value = 2743.8572
$ws.Range("L74").Value = 7833.222
$ws.Range("M74").Value = -1869.8572
$ws.Range("N74").Value = -9581.222
$ws.Range("H77").Value = 3981.8108
$ws.Range("I77").Value = 2743.8572
$ws.Range("J77").Value = 7833.222
$ws.Range("K77").Value = 13719.286
$ws.Range("L77").Value = 39166.11
$ws.Range("M77").Value = -9351.286
$ws.Range("N77").Value = -47902.11
$ws.Range("H122").Value = 2864.138
$ws.Range("I122").Value = 2505.3157
$ws.Range("K122").Value = 7515.9471
$ws.Range("M122").Value = -5065.9471
$ws.Range("H136").Value = 6267.2856
$ws.Range("I136").Value = 5313.121
$ws.Range("J136").Value = 9765.888999999999
$ws.Range("K136").Value = 15939.363
$ws.Range("L136").Value = 29297.667
$ws.Range("M136").Value = -13389.363
$ws.Range("N136").Value = -34397.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1471.0555
$ws.Range("I86").Value = 1471.0555
$ws.Range("K86").Value = 1471.0555
$ws.Range("M86").Value = -348.0554999999999
$ws.Range("H89").Value = 1471.0555
$ws.Range("I89").Value = 1471.0555
$ws.Range("K89").Value = 7355.2775
$ws.Range("M89").Value = -1739.2775
$ws.Range("H105").Value = 311.25
$ws.Range("I105").Value = 311.25
$ws.Range("K105").Value = 311.25
$ws.Range("M105").Value = 1435.75
$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 4000
$ws.Range("K107").Value = 4000
$ws.Range("M107").Value = -2080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 348.5
$ws.Range("I33").Value = 348.5
$ws.Range("K33").Value = 348.5
$ws.Range("M33").Value = 30.5
$ws.Range("H58").Value = 5413.0586
$ws.Range("I58").Value = 3440.1538
$ws.Range("K58").Value = 3440.1538
$ws.Range("M58").Value = -3237.1538
$ws.Range("H80").Value = 36788.99
$ws.Range("J80").Value = 36788.99
$ws.Range("L80").Value = 36788.99
$ws.Range("N80").Value = -39034.99
$ws.Range("H83").Value = 36788.99
$ws.Range("J83").Value = 36788.99
$ws.Range("L83").Value = 110366.97
$ws.Range("N83").Value = -121598.97
$ws.Range("H105").Value = 2943.7778
$ws.Range("I105").Value = 2670.8572
$ws.Range("K105").Value = 2670.8572
$ws.Range("M105").Value = -923.8571999999999
$ws.Range("H122").Value = 5211.2856
$ws.Range("I122").Value = 4685.9
$ws.Range("J122").Value = 6524.75
$ws.Range("K122").Value = 14057.7
$ws.Range("L122").Value = 19574.25
$ws.Range("M122").Value = -11607.7
$ws.Range("N122").Value = -24474.25
$ws.Range("H132").Value = 3378.9312
$ws.Range("I132").Value = 2987.4348
$ws.Range("K132").Value = 8962.304400000001
$ws.Range("M132").Value = -6432.304400000001
$ws.Range("H134").Value = 7096.476
$ws.Range("I134").Value = 6063.25
$ws.Range("K134").Value = 18189.75
$ws.Range("M134").Value = -15654.75
$ws.Range("H136").Value = 5413.0586
$ws.Range("I136").Value = 3440.1538
$ws.Range("K136").Value = 10320.4614
$ws.Range("M136").Value = -7770.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12540.5
$ws.Range("J2").Value = 20024.8
$ws.Range("L2").Value = 120148.8
$ws.Range("N2").Value = -120374.8
$ws.Range("H57").Value = 6000
$ws.Range("J57").Value = 6000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19118
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H107").Value = 379.25
$ws.Range("J107").Value = 379.25
$ws.Range("L107").Value = 1137.75
$ws.Range("N107").Value = -4977.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 11500
$ws.Range("I18").Value = 11500
$ws.Range("K18").Value = 11500
$ws.Range("M18").Value = -11207
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H70").Value = 15272.5
$ws.Range("I70").Value = 12249.75
$ws.Range("K70").Value = 12249.75
$ws.Range("M70").Value = -11979.75
$ws.Range("H73").Value = 15272.5
$ws.Range("I73").Value = 12249.75
$ws.Range("K73").Value = 12249.75
$ws.Range("M73").Value = -11313.75
$ws.Range("H80").Value = 6458.077
$ws.Range("I80").Value = 5942.8335
$ws.Range("K80").Value = 5942.8335
$ws.Range("M80").Value = -4944.8335
$ws.Range("H83").Value = 6458.077
$ws.Range("I83").Value = 5942.8335
$ws.Range("K83").Value = 29714.1675
$ws.Range("M83").Value = -24722.1675
$ws.Range("H96").Value = 33351.4
$ws.Range("J96").Value = 33351.4
$ws.Range("L96").Value = 33351.4
$ws.Range("N96").Value = -38843.4
$ws.Range("H100").Value = 20337.5
$ws.Range("J100").Value = 20337.5
$ws.Range("L100").Value = 20337.5
$ws.Range("N100").Value = -22501.5
$ws.Range("H111").Value = 84997.5
$ws.Range("J111").Value = 84997.5
$ws.Range("L111").Value = 84997.5
$ws.Range("N111").Value = -91131.5
$ws.Range("H132").Value = 2362
$ws.Range("I132").Value = 2074.4666
$ws.Range("J132").Value = 3799.6667
$ws.Range("K132").Value = 6223.399800000001
$ws.Range("L132").Value = 11399.0001
$ws.Range("M132").Value = -3693.399800000001
$ws.Range("N132").Value = -16459.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 14453.704
$ws.Range("I46").Value = 6274.7144
$ws.Range("J46").Value = 17316.35
$ws.Range("K46").Value = 6274.7144
$ws.Range("L46").Value = 17316.35
$ws.Range("M46").Value = -6086.7144
$ws.Range("N46").Value = -17692.35
$ws.Range("H93").Value = 2818.5
$ws.Range("I93").Value = 3144
$ws.Range("K93").Value = 3144
$ws.Range("M93").Value = -1896
$ws.Range("H132").Value = 8934.477000000001
$ws.Range("I132").Value = 9277.852999999999
$ws.Range("K132").Value = 27833.559
$ws.Range("M132").Value = -25303.559

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3486.513
$ws.Range("J122").Value = 4411.3335
$ws.Range("L122").Value = 13234.0005
$ws.Range("N122").Value = -18134.0005
$ws.Range("H132").Value = 6131.7646
$ws.Range("I132").Value = 5853.5835
$ws.Range("J132").Value = 6799.4
$ws.Range("K132").Value = 17560.7505
$ws.Range("L132").Value = 20398.2
$ws.Range("M132").Value = -15030.7505
$ws.Range("N132").Value = -25458.2
